$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 97.81044791763058
$ws.Range("C7").Value = 97.83095069353061
$ws.Range("D7").Value = 97.8086549863321
$ws.Range("E7").Value = 97.83385868602366

$ws.Range("B8").Value = 97.40134127420878
$ws.Range("C8").Value = 97.22221848585441
$ws.Range("D8").Value = 97.31014111981189
$ws.Range("E8").Value = 97.28598360610086

$ws.Range("B9").Value = 95.9465541823639
$ws.Range("C9").Value = 95.96405177831045
$ws.Range("D9").Value = 95.93426442441375
$ws.Range("E9").Value = 96.03245963294147
